# Weekly data refresh: insert a new daily observation at row 16 ("Fecha"
# 2021-11-19), pushing the existing rows 16-51 down to 17-52 (row 51's
# data ends up at row 52). This mirrors the upstream commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 16; Excel shifts rows 16..51
# down to 17..52 and extends the used range to A1:R52 automatically.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new observation. Columns
# that simply repeat the prior row's static attributes (Mercado ID,
# Mercado, Región, Codreg, Categoría ID, Categoría, Variedad, Calidad,
# Unidad de comercialización, Origen, Kg o Unidades, Clasificación) are
# filled in alongside the changed measurement columns (Fecha, Volumen,
# Precio mínimo/máximo/promedio ponderado, Precio $/Kg).
$ws.Cells.Item(16, 1).Value = 11
$ws.Cells.Item(16, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(16, 3).Value = "Bíobío"
$ws.Cells.Item(16, 4).Value = 44519
$ws.Cells.Item(16, 5).Value = 8
$ws.Cells.Item(16, 6).Value = 100112001
$ws.Cells.Item(16, 7).Value = "Berenjena"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 220
$ws.Cells.Item(16, 11).Value = 7500
$ws.Cells.Item(16, 12).Value = 8000
$ws.Cells.Item(16, 13).Value = 7773
$ws.Cells.Item(16, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(16, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(16, 16).Value = 130
$ws.Cells.Item(16, 17).Value = 60
$ws.Cells.Item(16, 18).Value = "Hortaliza"
